$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$ws.Range("B22").Value  = "A dildo as quick as the wind`nSex toys HP consumption -10%"
$ws.Range("B149").Value = "Legendary dildo`nSex toys, HP consumption 0, 2 attacks"
$ws.Range("B152").Value = "Telescopic rod`nSex toys HP consumption -10%"
$ws.Range("B181").Value = "Vibe that evaporates water with super vibration`nSex toy water attribute special effect consumption HP -50%"
$ws.Range("B224").Value = "A dildo that imitates the tail of an incubus`nSex toys automatic recovery of energy +1 automatic recovery of HP +1"
$ws.Range("B266").Value = "It's also a problem to be too late`nRubber HP consumption -10%"
$ws.Range("B272").Value = "Medicinal rubber that relieves tiredness and stiffness`nRubber automatic recovery of HP +2"
$ws.Range("B274").Value = "Rubber that refreshes the mind and body with a nice scent`nRubber automatic recovery +2 HP automatic recovery +2"
$ws.Range("B276").Value = "Rubber that seems to improve blood flow`nRubber automatic recovery +2 HP automatic recovery +2"
$ws.Range("B278").Value = "Rubber that is good for your body and keeps you healthy every day`nRubber automatic recovery of HP +4"
$ws.Range("B413").Value = "Fairy bracelet, a ring for humans`nRing fascination invalid, aphrodisiac poison invalid, HP consumption -20%"
$ws.Range("B422").Value = "A ring that expresses calmness`nRing Excitement disabled Automatic recovery of HP +3"
$ws.Range("B435").Value = "A ring that gives you strength`nRing, weakness invalid, restraint invalid, HP automatic recovery +1"
$ws.Range("B443").Value = "Black cat shop special ring limited to customers`nRing action additional 5% maximum energy +200 maximum HP +20"
$ws.Range("B457").Value = "A ring with the magical power of a famous magician`nRing: Automatic recovery of HP +5"
$ws.Range("B565").Value = "Earrings that shine on the night of the crescent moon`nAccessory Automatic recovery of HP +1"
$ws.Range("B567").Value = "Earrings that shine on a half-moon night`nAccessory Automatic recovery of HP +2"
$ws.Range("B570").Value = "Earrings that shine on a full moon night`nAccessory Automatic recovery of HP +3"
$ws.Range("B630").Value = "Shina's pants"
